# Apply the "evolution en %" rework of the physical metrics table:
# - Rename column D header to the new % based label
# - Recompute column D as the % difference between Top 5 and Bottom 15 averages,
#   rounded to 2 decimals
# - Re-sort the metric rows by descending absolute value of that % difference
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 4).Value = "Diff. Top 5 avec Bottom 15 en %"

# Row 2: highdecel_count_full_tip
$ws.Cells.Item(2, 1).Value = "highdecel_count_full_tip"
$ws.Cells.Item(2, 2).Value = 43.07795953848586
$ws.Cells.Item(2, 3).Value = 40.18806295821776
$ws.Cells.Item(2, 4).Value = 7.19
$ws.Cells.Item(2, 5).Value = 1.804049587967636
$ws.Cells.Item(2, 6).Value = 3.418269469976772
$ws.Cells.Item(2, 7).Value = 40.02631578947368
$ws.Cells.Item(2, 8).Value = 35.70588235294117
$ws.Cells.Item(2, 9).Value = 44.67567567567568
$ws.Cells.Item(2, 10).Value = 47.56756756756756

# Row 3: sprint_distance_full_tip
$ws.Cells.Item(3, 1).Value = "sprint_distance_full_tip"
$ws.Cells.Item(3, 2).Value = 835.3134660976766
$ws.Cells.Item(3, 3).Value = 788.8824798265357
$ws.Cells.Item(3, 4).Value = 5.89
$ws.Cells.Item(3, 5).Value = 51.89725678631734
$ws.Cells.Item(3, 6).Value = 76.13923929296327
$ws.Cells.Item(3, 7).Value = 755.6052631578947
$ws.Cells.Item(3, 8).Value = 657.6111111111111
$ws.Cells.Item(3, 9).Value = 896.8333333333334
$ws.Cells.Item(3, 10).Value = 907.921052631579

# Row 4: total_distance_full_otip
$ws.Cells.Item(4, 1).Value = "total_distance_full_otip"
$ws.Cells.Item(4, 2).Value = 37976.47255413308
$ws.Cells.Item(4, 3).Value = 40250.32098466204
$ws.Cells.Item(4, 4).Value = -5.65
$ws.Cells.Item(4, 5).Value = 1218.967684573305
$ws.Cells.Item(4, 6).Value = 2968.429064875285
$ws.Cells.Item(4, 7).Value = 36469.05263157895
$ws.Cells.Item(4, 8).Value = 33863.13513513513
$ws.Cells.Item(4, 9).Value = 39737.52777777778
$ws.Cells.Item(4, 10).Value = 44247.27777777778

# Row 5: highaccel_count_full_tip
$ws.Cells.Item(5, 1).Value = "highaccel_count_full_tip"
$ws.Cells.Item(5, 2).Value = 23.06278647068121
$ws.Cells.Item(5, 3).Value = 21.87721370721371
$ws.Cells.Item(5, 4).Value = 5.42
$ws.Cells.Item(5, 5).Value = 2.246411429220373
$ws.Cells.Item(5, 6).Value = 2.629106640832466
$ws.Cells.Item(5, 7).Value = 19.81578947368421
$ws.Cells.Item(5, 8).Value = 18.32432432432432
$ws.Cells.Item(5, 9).Value = 25.69444444444444
$ws.Cells.Item(5, 10).Value = 27.89189189189189

# Row 6: meddecel_count_full_tip
$ws.Cells.Item(6, 1).Value = "meddecel_count_full_tip"
$ws.Cells.Item(6, 2).Value = 278.8151019440493
$ws.Cells.Item(6, 3).Value = 265.5865184592739
$ws.Cells.Item(6, 4).Value = 4.98
$ws.Cells.Item(6, 5).Value = 19.22662352853307
$ws.Cells.Item(6, 6).Value = 19.54720228188021
$ws.Cells.Item(6, 7).Value = 248.2105263157895
$ws.Cells.Item(6, 8).Value = 235.7647058823529
$ws.Cells.Item(6, 9).Value = 298.2432432432432
$ws.Cells.Item(6, 10).Value = 306.8648648648648

# Row 7: meddecel_count_full_otip
$ws.Cells.Item(7, 1).Value = "meddecel_count_full_otip"
$ws.Cells.Item(7, 2).Value = 315.8095463884938
$ws.Cells.Item(7, 3).Value = 330.1202995456557
$ws.Cells.Item(7, 4).Value = -4.34
$ws.Cells.Item(7, 5).Value = 22.34108746755974
$ws.Cells.Item(7, 6).Value = 27.9827234374394
$ws.Cells.Item(7, 7).Value = 290.7105263157895
$ws.Cells.Item(7, 8).Value = 267.6
$ws.Cells.Item(7, 9).Value = 350.5555555555555
$ws.Cells.Item(7, 10).Value = 373.8181818181818

# Row 8: total_distance_full_tip
$ws.Cells.Item(8, 1).Value = "total_distance_full_tip"
$ws.Cells.Item(8, 2).Value = 37912.29730519994
$ws.Cells.Item(8, 3).Value = 36368.98142873967
$ws.Cells.Item(8, 4).Value = 4.24
$ws.Cells.Item(8, 5).Value = 2066.412088921724
$ws.Cells.Item(8, 6).Value = 2328.777700280471
$ws.Cells.Item(8, 7).Value = 34637.57894736842
$ws.Cells.Item(8, 8).Value = 32626.35294117647
$ws.Cells.Item(8, 9).Value = 39666.86486486487
$ws.Cells.Item(8, 10).Value = 40923.35135135135

# Row 9: sprint_count_full_tip
$ws.Cells.Item(9, 1).Value = "sprint_count_full_tip"
$ws.Cells.Item(9, 2).Value = 42.03405247352616
$ws.Cells.Item(9, 3).Value = 40.33372418898735
$ws.Cells.Item(9, 4).Value = 4.22
$ws.Cells.Item(9, 5).Value = 2.570313890578731
$ws.Cells.Item(9, 6).Value = 3.89549543514882
$ws.Cells.Item(9, 7).Value = 38.31578947368421
$ws.Cells.Item(9, 8).Value = 35
$ws.Cells.Item(9, 9).Value = 45.07894736842105
$ws.Cells.Item(9, 10).Value = 47.02702702702702

# Row 10: medaccel_count_full_tip
$ws.Cells.Item(10, 1).Value = "medaccel_count_full_tip"
$ws.Cells.Item(10, 2).Value = 392.7206495969654
$ws.Cells.Item(10, 3).Value = 378.1705047880899
$ws.Cells.Item(10, 4).Value = 3.85
$ws.Cells.Item(10, 5).Value = 30.18429608135872
$ws.Cells.Item(10, 6).Value = 24.29285216229298
$ws.Cells.Item(10, 7).Value = 345.3684210526316
$ws.Cells.Item(10, 8).Value = 341.0882352941176
$ws.Cells.Item(10, 9).Value = 421.6756756756757
$ws.Cells.Item(10, 10).Value = 426.7837837837838

# Row 11: running_distance_full_tip
$ws.Cells.Item(11, 1).Value = "running_distance_full_tip"
$ws.Cells.Item(11, 2).Value = 6039.2
$ws.Cells.Item(11, 3).Value = 5823.339042275358
$ws.Cells.Item(11, 4).Value = 3.71
$ws.Cells.Item(11, 5).Value = 442.3750533663978
$ws.Cells.Item(11, 6).Value = 412.9608627670697
$ws.Cells.Item(11, 7).Value = 5611.783783783784
$ws.Cells.Item(11, 8).Value = 5104.131578947368
$ws.Cells.Item(11, 9).Value = 6685.394736842105
$ws.Cells.Item(11, 10).Value = 6560.72972972973

# Row 12: hi_distance_full_tip
$ws.Cells.Item(12, 1).Value = "hi_distance_full_tip"
$ws.Cells.Item(12, 2).Value = 3187.744705231547
$ws.Cells.Item(12, 3).Value = 3090.024047231663
$ws.Cells.Item(12, 4).Value = 3.16
$ws.Cells.Item(12, 5).Value = 215.7537026835957
$ws.Cells.Item(12, 6).Value = 241.9822915721011
$ws.Cells.Item(12, 7).Value = 2908.28947368421
$ws.Cells.Item(12, 8).Value = 2756.631578947368
$ws.Cells.Item(12, 9).Value = 3415.578947368421
$ws.Cells.Item(12, 10).Value = 3483.297297297298

# Row 13: hi_count_full_tip
$ws.Cells.Item(13, 1).Value = "hi_count_full_tip"
$ws.Cells.Item(13, 2).Value = 265.7587719298245
$ws.Cells.Item(13, 3).Value = 257.9934589215085
$ws.Cells.Item(13, 4).Value = 3.01
$ws.Cells.Item(13, 5).Value = 20.0632245474992
$ws.Cells.Item(13, 6).Value = 19.25306162275538
$ws.Cells.Item(13, 7).Value = 238.8947368421053
$ws.Cells.Item(13, 8).Value = 227.7058823529412
$ws.Cells.Item(13, 9).Value = 286.8157894736842
$ws.Cells.Item(13, 10).Value = 292.4864864864865

# Row 14: running_distance_full_otip
$ws.Cells.Item(14, 1).Value = "running_distance_full_otip"
$ws.Cells.Item(14, 2).Value = 7319.500537379485
$ws.Cells.Item(14, 3).Value = 7536.131729121264
$ws.Cells.Item(14, 4).Value = -2.87
$ws.Cells.Item(14, 5).Value = 269.5763732664105
$ws.Cells.Item(14, 6).Value = 866.6375494370427
$ws.Cells.Item(14, 7).Value = 6919.815789473684
$ws.Cells.Item(14, 8).Value = 6206.216216216216
$ws.Cells.Item(14, 9).Value = 7561.432432432433
$ws.Cells.Item(14, 10).Value = 9735.833333333334

# Row 15: hsr_count_full_tip
$ws.Cells.Item(15, 1).Value = "hsr_count_full_tip"
$ws.Cells.Item(15, 2).Value = 223.7247194562984
$ws.Cells.Item(15, 3).Value = 217.6597347325211
$ws.Cells.Item(15, 4).Value = 2.79
$ws.Cells.Item(15, 5).Value = 17.55049506098603
$ws.Cells.Item(15, 6).Value = 15.53115578420854
$ws.Cells.Item(15, 7).Value = 200.578947368421
$ws.Cells.Item(15, 8).Value = 192.7058823529412
$ws.Cells.Item(15, 9).Value = 241.7368421052632
$ws.Cells.Item(15, 10).Value = 245.4594594594595

# Row 16: highaccel_count_full_otip
$ws.Cells.Item(16, 1).Value = "highaccel_count_full_otip"
$ws.Cells.Item(16, 2).Value = 22.41780464675201
$ws.Cells.Item(16, 3).Value = 21.84378565892498
$ws.Cells.Item(16, 4).Value = 2.63
$ws.Cells.Item(16, 5).Value = 1.766957316097681
$ws.Cells.Item(16, 6).Value = 1.142454331424871
$ws.Cells.Item(16, 7).Value = 20.2972972972973
$ws.Cells.Item(16, 8).Value = 20.23529411764706
$ws.Cells.Item(16, 9).Value = 25.08333333333333
$ws.Cells.Item(16, 10).Value = 24

# Row 17: medaccel_count_full_otip
$ws.Cells.Item(17, 1).Value = "medaccel_count_full_otip"
$ws.Cells.Item(17, 2).Value = 441.4351351351351
$ws.Cells.Item(17, 3).Value = 451.7737973875126
$ws.Cells.Item(17, 4).Value = -2.29
$ws.Cells.Item(17, 5).Value = 20.47642321843283
$ws.Cells.Item(17, 6).Value = 31.76771400089741
$ws.Cells.Item(17, 7).Value = 423.1351351351352
$ws.Cells.Item(17, 8).Value = 385.0857142857143
$ws.Cells.Item(17, 9).Value = 474
$ws.Cells.Item(17, 10).Value = 494.75

# Row 18: hsr_distance_full_tip
$ws.Cells.Item(18, 1).Value = "hsr_distance_full_tip"
$ws.Cells.Item(18, 2).Value = 2352.431239133871
$ws.Cells.Item(18, 3).Value = 2301.141567405127
$ws.Cells.Item(18, 4).Value = 2.23
$ws.Cells.Item(18, 5).Value = 180.0706048595034
$ws.Cells.Item(18, 6).Value = 174.9257849494258
$ws.Cells.Item(18, 7).Value = 2152.684210526316
$ws.Cells.Item(18, 8).Value = 2028.736842105263
$ws.Cells.Item(18, 9).Value = 2576.342105263158
$ws.Cells.Item(18, 10).Value = 2595.108108108108

# Row 19: sprint_count_full_otip
$ws.Cells.Item(19, 1).Value = "sprint_count_full_otip"
$ws.Cells.Item(19, 2).Value = 39.68480322427691
$ws.Cells.Item(19, 3).Value = 40.57962779433367
$ws.Cells.Item(19, 4).Value = -2.21
$ws.Cells.Item(19, 5).Value = 2.020522280728649
$ws.Cells.Item(19, 6).Value = 7.595268047957474
$ws.Cells.Item(19, 7).Value = 36.64864864864865
$ws.Cells.Item(19, 8).Value = 31.75675675675676
$ws.Cells.Item(19, 9).Value = 42.2972972972973
$ws.Cells.Item(19, 10).Value = 63.97222222222222

# Row 20: sprint_distance_full_otip
$ws.Cells.Item(20, 1).Value = "sprint_distance_full_otip"
$ws.Cells.Item(20, 2).Value = 755.829619092777
$ws.Cells.Item(20, 3).Value = 772.7644975038164
$ws.Cells.Item(20, 4).Value = -2.19
$ws.Cells.Item(20, 5).Value = 41.01112008095969
$ws.Cells.Item(20, 6).Value = 132.2968363426853
$ws.Cells.Item(20, 7).Value = 698.6756756756756
$ws.Cells.Item(20, 8).Value = 597.5526315789474
$ws.Cells.Item(20, 9).Value = 812.7297297297297
$ws.Cells.Item(20, 10).Value = 1165.694444444444

# Row 21: hi_count_full_otip
$ws.Cells.Item(21, 1).Value = "hi_count_full_otip"
$ws.Cells.Item(21, 2).Value = 304.9615457562826
$ws.Cells.Item(21, 3).Value = 308.5424472629117
$ws.Cells.Item(21, 4).Value = -1.16
$ws.Cells.Item(21, 5).Value = 16.16952708317574
$ws.Cells.Item(21, 6).Value = 43.62314153208373
$ws.Cells.Item(21, 7).Value = 279.8648648648648
$ws.Cells.Item(21, 8).Value = 250.1621621621622
$ws.Cells.Item(21, 9).Value = 324.1621621621622
$ws.Cells.Item(21, 10).Value = 437.4444444444445

# Row 22: hsr_count_full_otip
$ws.Cells.Item(22, 1).Value = "hsr_count_full_otip"
$ws.Cells.Item(22, 2).Value = 265.2767425320056
$ws.Cells.Item(22, 3).Value = 267.962819468578
$ws.Cells.Item(22, 4).Value = -1
$ws.Cells.Item(22, 5).Value = 14.19774683472161
$ws.Cells.Item(22, 6).Value = 36.15804972548949
$ws.Cells.Item(22, 7).Value = 243.2162162162162
$ws.Cells.Item(22, 8).Value = 218.4054054054054
$ws.Cells.Item(22, 9).Value = 281.8648648648648
$ws.Cells.Item(22, 10).Value = 373.4722222222222

# Row 23: highdecel_count_full_otip
$ws.Cells.Item(23, 1).Value = "highdecel_count_full_otip"
$ws.Cells.Item(23, 2).Value = 55.24603287498024
$ws.Cells.Item(23, 3).Value = 55.60428005545653
$ws.Cells.Item(23, 4).Value = -0.64
$ws.Cells.Item(23, 5).Value = 4.373797607140837
$ws.Cells.Item(23, 6).Value = 6.117368018980944
$ws.Cells.Item(23, 7).Value = 50.7027027027027
$ws.Cells.Item(23, 8).Value = 47.10810810810811
$ws.Cells.Item(23, 9).Value = 60.30555555555556
$ws.Cells.Item(23, 10).Value = 71.22222222222223

# Row 24: hsr_distance_full_otip
$ws.Cells.Item(24, 1).Value = "hsr_distance_full_otip"
$ws.Cells.Item(24, 2).Value = 2889.921906116643
$ws.Cells.Item(24, 3).Value = 2878.475614489082
$ws.Cells.Item(24, 4).Value = 0.4
$ws.Cells.Item(24, 5).Value = 151.9073601325664
$ws.Cells.Item(24, 6).Value = 455.2798591713864
$ws.Cells.Item(24, 7).Value = 2683.810810810811
$ws.Cells.Item(24, 8).Value = 2372.567567567567
$ws.Cells.Item(24, 9).Value = 3112.351351351351
$ws.Cells.Item(24, 10).Value = 4224.083333333333

# Row 25: total_metersperminute_full_tip
$ws.Cells.Item(25, 1).Value = "total_metersperminute_full_tip"
$ws.Cells.Item(25, 2).Value = 1931.688721510985
$ws.Cells.Item(25, 3).Value = 1926.662618168214
$ws.Cells.Item(25, 4).Value = 0.26
$ws.Cells.Item(25, 5).Value = 124.2866861270833
$ws.Cells.Item(25, 6).Value = 47.7370640042684
$ws.Cells.Item(25, 7).Value = 1802.196111111111
$ws.Cells.Item(25, 8).Value = 1859.296
$ws.Cells.Item(25, 9).Value = 2112.557894736842
$ws.Cells.Item(25, 10).Value = 2017.006756756757

# Row 26: hi_distance_full_otip
$ws.Cells.Item(26, 1).Value = "hi_distance_full_otip"
$ws.Cells.Item(26, 2).Value = 3645.75152520942
$ws.Cells.Item(26, 3).Value = 3651.240111992898
$ws.Cells.Item(26, 4).Value = -0.15
$ws.Cells.Item(26, 5).Value = 192.1082580171873
$ws.Cells.Item(26, 6).Value = 585.3472947855797
$ws.Cells.Item(26, 7).Value = 3382.486486486487
$ws.Cells.Item(26, 8).Value = 2993.783783783784
$ws.Cells.Item(26, 9).Value = 3925.081081081081
$ws.Cells.Item(26, 10).Value = 5389.777777777777

# Row 27: total_metersperminute_full_otip
$ws.Cells.Item(27, 1).Value = "total_metersperminute_full_otip"
$ws.Cells.Item(27, 2).Value = 2052.954703255888
$ws.Cells.Item(27, 3).Value = 2054.780069440093
$ws.Cells.Item(27, 4).Value = -0.09
$ws.Cells.Item(27, 5).Value = 132.5511653665204
$ws.Cells.Item(27, 6).Value = 82.6271374705141
$ws.Cells.Item(27, 7).Value = 1855.891944444445
$ws.Cells.Item(27, 8).Value = 1919.415405405405
$ws.Cells.Item(27, 9).Value = 2207.652162162162
$ws.Cells.Item(27, 10).Value = 2272.552222222223
